$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16:
# A: Class, B: harmonized_filename, C: image_path_blur, D: image_path_no_blur, E: Logistic_Regression_Classification
$classes = @("Sedan","Sedan","Pickup","SUV","Sedan","Convertible","Pickup","SUV","Pickup","Sedan","SUV","SUV","Pickup","Pickup","Sedan")
$files   = @("Sedan_test_orig_test_00697_resized.jpg","Sedan_test_orig_test_04302_resized.jpg","Pickup_test_orig_test_07661_resized.jpg","SUV_test_orig_test_02156_resized.jpg","Sedan_test_orig_test_01820_resized.jpg","Convertible_test_orig_test_06347_resized.jpg","Pickup_test_orig_test_00858_resized.jpg","SUV_test_orig_test_03452_resized.jpg","Pickup_test_orig_train_02871_resized.jpg","Sedan_test_orig_test_06356_resized.jpg","SUV_test_orig_test_01769_resized.jpg","SUV_test_orig_train_03247_resized.jpg","Pickup_test_orig_test_04440_resized.jpg","SUV_test_orig_train_04481_resized.jpg","Sedan_test_orig_train_08005_resized.jpg")
$labels  = @("Sedan","Sedan","SUV","Convertible","Sedan","Convertible","Pickup","Sedan","Pickup","Sedan","Sedan","Sedan","SUV","SUV","Sedan")

for ($i = 0; $i -lt $classes.Length; $i++) {
    $rowIndex = $i + 2
    $filename = $files[$i]

    $ws.Cells.Item($rowIndex, 1).Value = $classes[$i]
    $ws.Cells.Item($rowIndex, 2).Value = $filename
    $ws.Cells.Item($rowIndex, 3).Value = "../../../Images/test/Blurred/" + $filename
    $ws.Cells.Item($rowIndex, 4).Value = "../../../Images/test/No_Blur/" + $filename
    $ws.Cells.Item($rowIndex, 5).Value = $labels[$i]
}
